$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) / "Valor Mora" (column F) data block occupies
# rows 16-46. The commit reverses the vertical order of that block (oldest
# periods were at the bottom, now they are at the top and vice-versa), while
# leaving every other column (B,C,D,G,H,I,J) untouched since their content is
# identical on every row.

$firstRow = 16
$lastRow = 46

# Snapshot the current Periodo Mora / Valor Mora pairs before overwriting
# anything, so the read side is not affected by the writes we are about to do.
$periodos = @()
$valores = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += $ws.Cells.Item($r, 5).Value()
    $valores += $ws.Cells.Item($r, 6).Value()
}

$count = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $count; $i++) {
    $srcIndex = $count - 1 - $i
    $destRow = $firstRow + $i
    $ws.Cells.Item($destRow, 5).Value = $periodos[$srcIndex]
    $ws.Cells.Item($destRow, 6).Value = $valores[$srcIndex]
}
